$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "48.390.23"
$ws.Range("E2").Value = "  +2.05%  "
$ws.Range("D3").Value = "2.529.06"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.527"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +3.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.25"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0821"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.63"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.29%  "
$ws.Range("E13").Value = "  +1.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("D15").Value = "2.924.16"
$ws.Range("E15").Value = "  +0.66%  "
$ws.Range("D16").Value = "2.541.18"
$ws.Range("E16").Value = "  +1.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.857"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("D18").Value = "48.278.06"
$ws.Range("E18").Value = "  +1.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.87%  "
$ws.Range("D21").Value = "0.0₃0947"
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("E28").Value = "  +1.52%  "
$ws.Range("E29").Value = "  +4.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.32"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -9.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.92"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.96"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("E34").Value = "  -1.09%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0791"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.72"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.99"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.36%  "
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("E41").Value = "  +4.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "119.42"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.19"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0301"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.63%  "
$ws.Range("D45").Value = "2.003.55"
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.13"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("E47").Value = "  -2.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.86"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.12"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("E50").Value = "  -0.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.37"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.17%  "
